# Updates the daily MeteoCat summary extraction:
# - advances DATA_DIA / DATA_EXTRACCIO / URL_FONT dates from 2026-02-14 to 2026-02-15
# - the 2026-02-15 00:xx run found no fresh station data, so every metric column
#   that previously held a reading is replaced with the "sense dades" placeholder
#   (two placeholders for columns that used to carry a value + a timestamp)
# - PRESSIO_ATMOSFERICA (column J) had stray figures baked in; those are cleared
# - widens columns L/M/N (RATXA_VENT_MAX / TEMPERATURA_MAXIMA_DIA / TEMPERATURA_MINIMA_DIA)
#   to a uniform width of 25, matching column O

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ r=2; D='2026-02-15'; E='2026-02-15 00:18:47'; F='https://www.meteo.cat/observacions/xema/dades?codi=YT&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L=$null; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=3; D='2026-02-15'; E='2026-02-15 00:18:49'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z1&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=4; D='2026-02-15'; E='2026-02-15 00:18:52'; F='https://www.meteo.cat/observacions/xema/dades?codi=DN&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=5; D='2026-02-15'; E='2026-02-15 00:18:55'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z6&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=6; D='2026-02-15'; E='2026-02-15 00:18:57'; F='https://www.meteo.cat/observacions/xema/dades?codi=DJ&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=7; D='2026-02-15'; E='2026-02-15 00:19:00'; F='https://www.meteo.cat/observacions/xema/dades?codi=X4&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=8; D='2026-02-15'; E='2026-02-15 00:19:02'; F='https://www.meteo.cat/observacions/xema/dades?codi=D5&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=9; D='2026-02-15'; E='2026-02-15 00:19:05'; F='https://www.meteo.cat/observacions/xema/dades?codi=YS&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=10; D='2026-02-15'; E='2026-02-15 00:19:07'; F='https://www.meteo.cat/observacions/xema/dades?codi=UN&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=11; D='2026-02-15'; E='2026-02-15 00:19:10'; F='https://www.meteo.cat/observacions/xema/dades?codi=MS&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K=$null; L=$null; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=12; D='2026-02-15'; E='2026-02-15 00:19:12'; F='https://www.meteo.cat/observacions/xema/dades?codi=W1&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K=$null; L=$null; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=13; D='2026-02-15'; E='2026-02-15 00:19:15'; F='https://www.meteo.cat/observacions/xema/dades?codi=DP&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=14; D='2026-02-15'; E='2026-02-15 00:19:17'; F='https://www.meteo.cat/observacions/xema/dades?codi=XL&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=15; D='2026-02-15'; E='2026-02-15 00:19:20'; F='https://www.meteo.cat/observacions/xema/dades?codi=VZ&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K=$null; L=$null; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=16; D='2026-02-15'; E='2026-02-15 00:19:22'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z7&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=17; D='2026-02-15'; E='2026-02-15 00:19:25'; F='https://www.meteo.cat/observacions/xema/dades?codi=XK&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=18; D='2026-02-15'; E='2026-02-15 00:19:28'; F='https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=19; D='2026-02-15'; E='2026-02-15 00:19:30'; F='https://www.meteo.cat/observacions/xema/dades?codi=YU&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=20; D='2026-02-15'; E='2026-02-15 00:19:33'; F='https://www.meteo.cat/observacions/xema/dades?codi=ZE&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=21; D='2026-02-15'; E='2026-02-15 00:19:35'; F='https://www.meteo.cat/observacions/xema/dades?codi=CD&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=22; D='2026-02-15'; E='2026-02-15 00:19:38'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z2&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=23; D='2026-02-15'; E='2026-02-15 00:19:40'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z5&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=24; D='2026-02-15'; E='2026-02-15 00:19:43'; F='https://www.meteo.cat/observacions/xema/dades?codi=VK&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=25; D='2026-02-15'; E='2026-02-15 00:19:45'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z3&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=26; D='2026-02-15'; E='2026-02-15 00:19:48'; F='https://www.meteo.cat/observacions/xema/dades?codi=CG&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=27; D='2026-02-15'; E='2026-02-15 00:19:50'; F='https://www.meteo.cat/observacions/xema/dades?codi=Z9&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=28; D='2026-02-15'; E='2026-02-15 00:19:53'; F='https://www.meteo.cat/observacions/xema/dades?codi=YB&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=29; D='2026-02-15'; E='2026-02-15 00:19:55'; F='https://www.meteo.cat/observacions/xema/dades?codi=YP&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=30; D='2026-02-15'; E='2026-02-15 00:19:58'; F='https://www.meteo.cat/observacions/xema/dades?codi=J5&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=31; D='2026-02-15'; E='2026-02-15 00:20:00'; F='https://www.meteo.cat/observacions/xema/dades?codi=D6&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=32; D='2026-02-15'; E='2026-02-15 00:20:03'; F='https://www.meteo.cat/observacions/xema/dades?codi=XR&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=33; D='2026-02-15'; E='2026-02-15 00:20:05'; F='https://www.meteo.cat/observacions/xema/dades?codi=YA&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=34; D='2026-02-15'; E='2026-02-15 00:20:08'; F='https://www.meteo.cat/observacions/xema/dades?codi=DG&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=35; D='2026-02-15'; E='2026-02-15 00:20:10'; F='https://www.meteo.cat/observacions/xema/dades?codi=X5&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=36; D='2026-02-15'; E='2026-02-15 00:20:13'; F='https://www.meteo.cat/observacions/xema/dades?codi=D4&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=37; D='2026-02-15'; E='2026-02-15 00:20:15'; F='https://www.meteo.cat/observacions/xema/dades?codi=CI&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K=$null; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=38; D='2026-02-15'; E='2026-02-15 00:20:18'; F='https://www.meteo.cat/observacions/xema/dades?codi=XS&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=39; D='2026-02-15'; E='2026-02-15 00:20:20'; F='https://www.meteo.cat/observacions/xema/dades?codi=ZC&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=40; D='2026-02-15'; E='2026-02-15 00:20:22'; F='https://www.meteo.cat/observacions/xema/dades?codi=XH&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=''; K=$null; L=$null; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=41; D='2026-02-15'; E='2026-02-15 00:20:25'; F='https://www.meteo.cat/observacions/xema/dades?codi=XE&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=42; D='2026-02-15'; E='2026-02-15 00:20:27'; F='https://www.meteo.cat/observacions/xema/dades?codi=UE&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K=$null; L=$null; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=43; D='2026-02-15'; E='2026-02-15 00:20:30'; F='https://www.meteo.cat/observacions/xema/dades?codi=XO&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=44; D='2026-02-15'; E='2026-02-15 00:20:32'; F='https://www.meteo.cat/observacions/xema/dades?codi=VS&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=$null; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=45; D='2026-02-15'; E='2026-02-15 00:20:35'; F='https://www.meteo.cat/observacions/xema/dades?codi=YN&dia=2026-02-15T09:00Z'; G='sense dades'; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
    @{ r=46; D='2026-02-15'; E='2026-02-15 00:20:37'; F='https://www.meteo.cat/observacions/xema/dades?codi=D7&dia=2026-02-15T09:00Z'; G=$null; H='sense dades'; I='sense dades'; J=''; K='sense dades'; L='sense dades sense dades'; M='sense dades sense dades'; N='sense dades sense dades'; O='sense dades' },
)

foreach ($row in $rowData) {
    $r = $row.r

    # --- column D (DATA_DIA) -------------------------------------------------
    # Plain "yyyy-mm-dd" text gets auto-recognised as a date by the COM layer,
    # which would flip the cell to a numeric date serial + new number format.
    # Force it to text first, write it, then restore the original cell format
    # (border/alignment, no custom number format) via a formats-only paste from
    # a neighbouring cell that already carries that exact style.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row.D
    $ws.Cells.Item($r, 1).Copy() | Out-Null
    $dCell.PasteSpecial(-4122) | Out-Null

    # --- column E (DATA_EXTRACCIO) -------------------------------------------
    $ws.Cells.Item($r, 5).Value = $row.E

    # --- column F (URL_FONT) --------------------------------------------------
    $ws.Cells.Item($r, 6).Value = $row.F

    # --- column G (GRUIX_NEU_MAX) ----------------------------------------------
    if ($row.G -ne $null) { $ws.Cells.Item($r, 7).Value = $row.G }

    # --- column H (HUMITAT_MITJANA_DIA) ----------------------------------------
    if ($row.H -ne $null) { $ws.Cells.Item($r, 8).Value = $row.H }

    # --- column I (PRECIPITACIO_ACUM_DIA) ---------------------------------------
    if ($row.I -ne $null) { $ws.Cells.Item($r, 9).Value = $row.I }

    # --- column J (PRESSIO_ATMOSFERICA) -----------------------------------------
    if ($row.J -eq "") { $ws.Cells.Item($r, 10).ClearContents() }

    # --- column K (RADIACIO_GLOBAL) ---------------------------------------------
    if ($row.K -ne $null) { $ws.Cells.Item($r, 11).Value = $row.K }

    # --- column L (RATXA_VENT_MAX) -----------------------------------------------
    if ($row.L -ne $null) { $ws.Cells.Item($r, 12).Value = $row.L }

    # --- column M (TEMPERATURA_MAXIMA_DIA) ---------------------------------------
    if ($row.M -ne $null) { $ws.Cells.Item($r, 13).Value = $row.M }

    # --- column N (TEMPERATURA_MINIMA_DIA) ----------------------------------------
    if ($row.N -ne $null) { $ws.Cells.Item($r, 14).Value = $row.N }

    # --- column O (TEMPERATURA_MITJANA_DIA) -----------------------------------------
    if ($row.O -ne $null) { $ws.Cells.Item($r, 15).Value = $row.O }
}

# --- widen columns L, M, N to 25 characters (matches column O) -------------------
# ColumnWidth is expressed in the "characters of the default font" unit that Excel
# uses in its UI, which does not map 1:1 onto the raw OOXML column width; 24.17 is
# the value that round-trips to an OOXML width of exactly 25 for this workbook's
# default font/theme.
$ws.Columns("L").ColumnWidth = 24.17
$ws.Columns("M").ColumnWidth = 24.17
$ws.Columns("N").ColumnWidth = 24.17
